$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "US Core PMO ServiceRequest Profile" (row 43). This
# shifts all subsequent rows up by one and Excel automatically leaves the
# trailing row (57) blank/removed.
$ws.Rows("43:43").Delete()
